$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only cells in column D (Price) need the Text number format forced,
# since their values are numeric-looking strings (e.g. "300.28",
# "0.02690", "21.00") that would otherwise be auto-parsed into numbers
# and lose significant trailing/leading zeros. Column E (Volume) values
# contain "%" and spaces so Excel already keeps them as text.
# Cells that are not being modified are left completely untouched so
# their existing style/format is not disturbed.

foreach ($addr in @("D2","D3","D6","D7","D8","D9","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "23.437.25"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.643.65"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "300.28"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").Value = "0.3786"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("D8").Value = "50.44"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").Value = "0.3502"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "22.08"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").Value = "6.289"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").Value = "7.246"
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("D16").Value = "0.00001210"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "1.650.33"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "95.23"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D19").Value = "0.06991"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "6.630"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").Value = "17.39"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "12.42"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("D24").Value = "23.475.39"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "2.413"
$ws.Range("E25").Value = "  -4.23%  "
$ws.Range("D26").Value = "3.009"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "21.00"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").Value = "151.85"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "5.188"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").Value = "131.64"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("D31").Value = "1.830.62"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "6.845"
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("D33").Value = "2.145"
$ws.Range("E33").Value = "  -4.59%  "
$ws.Range("D34").Value = "11.16"
$ws.Range("E34").Value = "  -7.14%  "
$ws.Range("D35").Value = "0.9872"
$ws.Range("E35").Value = "  -7.07%  "
$ws.Range("D36").Value = "0.02690"
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("D37").Value = "0.08779"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "5.903"
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("D39").Value = "0.2416"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("D40").Value = "0.06777"
$ws.Range("E40").Value = "  -3.21%  "
$ws.Range("D41").Value = "12.85"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "0.6881"
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").Value = "1.293"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").Value = "15.54"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "0.6384"
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("D47").Value = "3.926"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "2.238"
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("D49").Value = "127.42"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "0.07665"
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("D51").Value = "1.237"
$ws.Range("E51").Value = "  +2.76%  "
